$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new blank column before column A.
#    This shifts the old column A -> B and old column C -> D,
#    carrying all values & formatting with it.
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()

# ---------------------------------------------------------------------------
# 2. Reorder rows 3,4,5 (3-cycle: old3->5, old4->3, old5->4) using a staging
#    row far below the used range so we never overwrite data we still need.
# ---------------------------------------------------------------------------
$ws.Range("B3:D3").Cut($ws.Range("B50:D50"))      # stash old row3 (the "however.. location" row)
$ws.Range("B4:D4").Cut($ws.Range("B3:D3"))        # old row4 ("Mismatch") -> row3
$ws.Range("B5:D5").Cut($ws.Range("B4:D4"))        # old row5 ("DLL does not appear") -> row4
$ws.Range("B50:D50").Cut($ws.Range("B5:D5"))      # stashed row3 ("however.. location") -> row5

# ---------------------------------------------------------------------------
# 3. Move the last row (old row7, "doesn't exist") down to row 11, leaving
#    rows 7-10 empty.
# ---------------------------------------------------------------------------
$ws.Range("B7:D7").Cut($ws.Range("B11:D11"))

# ---------------------------------------------------------------------------
# 4. Update the fill colours that changed (in place, by current content).
#    Row2 (exists,Class ID..): F5DEB3 -> 90FA90
#    Row3 (Mismatch, just moved here): FFFF00 -> E8FEE8
#    Row4 (DLL does not appear, just moved here): 98FB98 -> FFE699
# ---------------------------------------------------------------------------
$ws.Range("A2:B2").Interior.Color = 0x90FA90
$ws.Range("A3:B3").Interior.Color = 0xE8FEE8
$ws.Range("A4:B4").Interior.Color = 0xFFE699

# ---------------------------------------------------------------------------
# 5. Style column A cells per the target layout (column A is now empty,
#    only the fills/format matter).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5:A6").Interior.Color = 0xFF0000
$ws.Range("A11").Interior.Color = 0xFF0000
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Replace the D-column labels (they now show colour-name text rather than
#    matching the old C column order) and add the new RGB helper columns
#    F/G/H plus the new J11 note.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "00FF00"
$ws.Range("F1").Value = 0
$ws.Range("G1").Value = 255
$ws.Range("H1").Value = 0

$ws.Range("D2").Value = "BAFCBA"
$ws.Range("F2").Value = 186
$ws.Range("G2").Value = 252
$ws.Range("H2").Value = 186

$ws.Range("D3").Value = "E8FEE8"
$ws.Range("F3").Value = 232
$ws.Range("G3").Value = 254
$ws.Range("H3").Value = 232

$ws.Range("D4").Value = "FFE699"
$ws.Range("F4").Value = 255
$ws.Range("G4").Value = 230
$ws.Range("H4").Value = 153

$ws.Range("D5").Value = "FFC0CB"
$ws.Range("F5").Value = 255
$ws.Range("G5").Value = 192
$ws.Range("H5").Value = 203

$ws.Range("D6").Value = "FF69B4"
$ws.Range("F6").Value = 255
$ws.Range("G6").Value = 105
$ws.Range("H6").Value = 180

$ws.Range("D11").Value = "FF0000"
$ws.Range("F11").Value = 255
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = "Can't exist!"

# ---------------------------------------------------------------------------
# 7. Misc sheet-level tweaks from the diff: new column width target is on
#    column B now (handled automatically by the column insert above), row 1
#    is taller (title row) and the active selection/cell changes.
# ---------------------------------------------------------------------------
$ws.Rows("1").RowHeight = 31.5
$ws.Range("A4").Select()

Write-Host "done"
